# Fix Word document corruption issue by removing the embedded Merlin logo
# image and replacing it with a lightweight text placeholder run, per:
#   "Replaced embedded logo image with text placeholder to resolve Word
#    'unreadable content' error ... Enhanced template now uses text-based
#    Merlin branding instead of embedded PNG."

$d = $word.ActiveDocument

# Locate the paragraph that hosts the inline picture (the Merlin logo)
# by scanning paragraphs for one whose range owns an InlineShape, rather
# than relying on a hard-coded paragraph index.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    Write-Host "No inline image paragraph found; nothing to do."
} else {
    $r = $targetPara.Range

    # Replace the whole paragraph (which currently is just a <w:r><w:drawing>...)
    # with a single run carrying the italic grey "[MERLIN_LOGO_PLACEHOLDER]"
    # text, keeping the same right-aligned paragraph formatting.
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
           '<w:pPr><w:jc w:val="right"/></w:pPr>' +
           '<w:r>' +
           '<w:rPr>' +
           '<w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/>' +
           '<w:i/>' +
           '<w:iCs/>' +
           '<w:color w:val="6B7280"/>' +
           '<w:sz w:val="16"/>' +
           '<w:szCs w:val="16"/>' +
           '</w:rPr>' +
           '<w:t xml:space="preserve">[MERLIN_LOGO_PLACEHOLDER]</w:t>' +
           '</w:r>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'

    $r.InsertXML($xml) | Out-Null
    Write-Host "Replaced embedded Merlin logo image with text placeholder."
}
